$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expansion List")
$ws.Activate()

# Insert a new blank row above row 13 (shifts existing data rows 13-21 down to 14-22)
$ws.Rows(13).EntireRow.Insert()

# Select the newly inserted row to match the resulting UI state
$ws.Range("A13:XFD13").Select()
